$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 7

$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 4

$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 4

$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 0

$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 8
